$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.488.29'
$ws.Range('E2').Value = '  -0.14%  '
$ws.Range('D3').Value = '1.806.15'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.86'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('E6').Value = '  +3.38%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.76'
$ws.Range('E8').Value = '  +6.61%  '
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0964'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '2.067.42'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.56'
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.652'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '1.805.66'
$ws.Range('E15').Value = '  -0.58%  '
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '34.486.78'
$ws.Range('E17').Value = '  -0.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.60'
$ws.Range('E18').Value = '  +0.76%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.47'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('E20').Value = '  -1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.60'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.19'
$ws.Range('E24').Value = '  +4.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.99'
$ws.Range('E25').Value = '  -0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.98'
$ws.Range('E26').Value = '  +7.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.95'
$ws.Range('E27').Value = '  +1.56%  '
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('E34').Value = '  -1.55%  '
$ws.Range('D35').Value = '1.396.06'
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.673'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  -5.40%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0191'
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '83.07'
$ws.Range('E40').Value = '  -3.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.966'
$ws.Range('E41').Value = '  +1.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.83'
$ws.Range('E42').Value = '  -0.79%  '
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('E44').Value = '  +7.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.49'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.03'
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0499'
$ws.Range('E47').Value = '  -5.23%  '
$ws.Range('D48').Value = '1.969.03'
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '104.62'
$ws.Range('E49').Value = '  -0.92%  '
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('E51').Value = '  -3.15%  '
